$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit re-orders the 34 data rows (2-35): every row's varying fields
# (Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Origen,
# Precio $/Kg) are taken from a different source row, producing a full
# permutation of the data block. Columns that are identical on every row
# (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria, Unidad de
# comercializacion, Kg o Unidades, Clasificacion) are left untouched.

$cols = @("D", "H", "I", "J", "K", "L", "M", "O", "P")

# Snapshot the "before" values for every row so writes don't clobber reads.
$snapshot = @{}
for ($r = 2; $r -le 35; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# destination row -> source row (values copied from source into destination)
$map = @{}
$map[2] = 12
$map[3] = 32
$map[4] = 33
$map[5] = 22
$map[6] = 20
$map[7] = 30
$map[8] = 25
$map[9] = 24
$map[10] = 15
$map[11] = 19
$map[12] = 9
$map[13] = 31
$map[14] = 17
$map[15] = 21
$map[16] = 6
$map[17] = 13
$map[18] = 5
$map[19] = 28
$map[20] = 14
$map[21] = 3
$map[22] = 29
$map[23] = 11
$map[24] = 2
$map[25] = 26
$map[26] = 7
$map[27] = 10
$map[28] = 35
$map[29] = 34
$map[30] = 18
$map[31] = 16
$map[32] = 4
$map[33] = 23
$map[34] = 27
$map[35] = 8

for ($destRow = 2; $destRow -le 35; $destRow++) {
    $srcRow = $map[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcData[$c]
    }
}
